$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(28, 8).Value = 2720.375
$ws_ALC.Cells.Item(28, 9).Value = 2720.375
$ws_ALC.Cells.Item(28, 11).Value = 2720.375
$ws_ALC.Cells.Item(28, 13).Value = -2235.375
$ws_ALC.Cells.Item(32, 8).Value = 2312.5
$ws_ALC.Cells.Item(32, 9).Value = 625
$ws_ALC.Cells.Item(32, 10).Value = 4000
$ws_ALC.Cells.Item(32, 11).Value = 625
$ws_ALC.Cells.Item(32, 12).Value = 4000
$ws_ALC.Cells.Item(32, 13).Value = -299
$ws_ALC.Cells.Item(32, 14).Value = -4652
$ws_ALC.Cells.Item(39, 8).Value = 245.61539
$ws_ALC.Cells.Item(39, 9).Value = 275.3
$ws_ALC.Cells.Item(39, 10).Value = 146.66667
$ws_ALC.Cells.Item(39, 11).Value = 825.9000000000001
$ws_ALC.Cells.Item(39, 12).Value = 440.00001
$ws_ALC.Cells.Item(39, 13).Value = -529.9000000000001
$ws_ALC.Cells.Item(39, 14).Value = -1032.00001
$ws_ALC.Cells.Item(53, 8).Value = 1700.0769
$ws_ALC.Cells.Item(53, 9).Value = 1991
$ws_ALC.Cells.Item(53, 10).Value = 100
$ws_ALC.Cells.Item(53, 11).Value = 1991
$ws_ALC.Cells.Item(53, 12).Value = 100
$ws_ALC.Cells.Item(53, 13).Value = -1354
$ws_ALC.Cells.Item(53, 14).Value = -1374
$ws_ALC.Cells.Item(62, 8).Value = 12349566
$ws_ALC.Cells.Item(62, 9).Value = 13892886
$ws_ALC.Cells.Item(62, 10).Value = 3006
$ws_ALC.Cells.Item(62, 11).Value = 13892886
$ws_ALC.Cells.Item(62, 12).Value = 3006
$ws_ALC.Cells.Item(62, 13).Value = -13892262
$ws_ALC.Cells.Item(62, 14).Value = -4254
$ws_ALC.Cells.Item(64, 8).Value = 3907.2856
$ws_ALC.Cells.Item(64, 9).Value = 3860.25
$ws_ALC.Cells.Item(64, 11).Value = 3860.25
$ws_ALC.Cells.Item(64, 13).Value = -3612.25
$ws_ALC.Cells.Item(65, 8).Value = 12349566
$ws_ALC.Cells.Item(65, 9).Value = 13892886
$ws_ALC.Cells.Item(65, 10).Value = 3006
$ws_ALC.Cells.Item(65, 11).Value = 69464430
$ws_ALC.Cells.Item(65, 12).Value = 15030
$ws_ALC.Cells.Item(65, 13).Value = -69461310
$ws_ALC.Cells.Item(65, 14).Value = -21270
$ws_ALC.Cells.Item(67, 8).Value = 3907.2856
$ws_ALC.Cells.Item(67, 9).Value = 3860.25
$ws_ALC.Cells.Item(67, 11).Value = 3860.25
$ws_ALC.Cells.Item(67, 13).Value = -3002.25
$ws_ALC.Cells.Item(92, 8).Value = 2241.6155
$ws_ALC.Cells.Item(92, 9).Value = 2652.6667
$ws_ALC.Cells.Item(92, 10).Value = 1889.2858
$ws_ALC.Cells.Item(92, 11).Value = 2652.6667
$ws_ALC.Cells.Item(92, 12).Value = 1889.2858
$ws_ALC.Cells.Item(92, 13).Value = -1404.6667
$ws_ALC.Cells.Item(92, 14).Value = -4385.2858
$ws_ALC.Cells.Item(107, 8).Value = 2211.318
$ws_ALC.Cells.Item(107, 9).Value = 2064.5
$ws_ALC.Cells.Item(107, 10).Value = 2387.5
$ws_ALC.Cells.Item(107, 11).Value = 2064.5
$ws_ALC.Cells.Item(107, 12).Value = 2387.5
$ws_ALC.Cells.Item(107, 13).Value = -144.5
$ws_ALC.Cells.Item(107, 14).Value = -6227.5
$ws_ALC.Cells.Item(113, 8).Value = 3108
$ws_ALC.Cells.Item(113, 9).Value = 3090
$ws_ALC.Cells.Item(113, 10).Value = 3120
$ws_ALC.Cells.Item(113, 11).Value = 3090
$ws_ALC.Cells.Item(113, 12).Value = 3120
$ws_ALC.Cells.Item(113, 13).Value = 164
$ws_ALC.Cells.Item(113, 14).Value = -9628
$ws_ALC.Cells.Item(116, 8).Value = 2890.4
$ws_ALC.Cells.Item(116, 9).Value = 2213.6365
$ws_ALC.Cells.Item(116, 10).Value = 4751.5
$ws_ALC.Cells.Item(116, 11).Value = 2213.6365
$ws_ALC.Cells.Item(116, 12).Value = 4751.5
$ws_ALC.Cells.Item(116, 13).Value = 1228.3635
$ws_ALC.Cells.Item(116, 14).Value = -11635.5
$ws_ALC.Cells.Item(129, 8).Value = 630.4286
$ws_ALC.Cells.Item(129, 9).Value = 571.1667
$ws_ALC.Cells.Item(129, 11).Value = 1713.5001
$ws_ALC.Cells.Item(129, 13).Value = 3286.4999
$ws_ALC.Cells.Item(135, 8).Value = 55557404
$ws_ALC.Cells.Item(135, 9).Value = 662.7143
$ws_ALC.Cells.Item(135, 10).Value = 250006000
$ws_ALC.Cells.Item(135, 11).Value = 5964.428699999999
$ws_ALC.Cells.Item(135, 12).Value = 2250054000
$ws_ALC.Cells.Item(135, 13).Value = -3429.428699999999
$ws_ALC.Cells.Item(135, 14).Value = -2250059070
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(32, 8).Value = 5241.7114
$ws_ARM.Cells.Item(32, 9).Value = 5191.8
$ws_ARM.Cells.Item(32, 10).Value = 6489.5
$ws_ARM.Cells.Item(32, 11).Value = 5191.8
$ws_ARM.Cells.Item(32, 12).Value = 6489.5
$ws_ARM.Cells.Item(32, 13).Value = -4904.8
$ws_ARM.Cells.Item(32, 14).Value = -7063.5
$ws_ARM.Cells.Item(61, 8).Value = 47620076
$ws_ARM.Cells.Item(61, 9).Value = 55556444
$ws_ARM.Cells.Item(61, 10).Value = 1853.3334
$ws_ARM.Cells.Item(61, 11).Value = 55556444
$ws_ARM.Cells.Item(61, 12).Value = 1853.3334
$ws_ARM.Cells.Item(61, 13).Value = -55556232
$ws_ARM.Cells.Item(61, 14).Value = -2277.3334
$ws_ARM.Cells.Item(102, 8).Value = 27779134
$ws_ARM.Cells.Item(102, 10).Value = 1005.5
$ws_ARM.Cells.Item(102, 12).Value = 1005.5
$ws_ARM.Cells.Item(102, 14).Value = -4249.5
$ws_ARM.Cells.Item(110, 8).Value = 1083.5294
$ws_ARM.Cells.Item(110, 9).Value = 640.4666999999999
$ws_ARM.Cells.Item(110, 10).Value = 4406.5
$ws_ARM.Cells.Item(110, 11).Value = 640.4666999999999
$ws_ARM.Cells.Item(110, 12).Value = 4406.5
$ws_ARM.Cells.Item(110, 13).Value = 1404.5333
$ws_ARM.Cells.Item(110, 14).Value = -8496.5
$ws_ARM.Cells.Item(122, 8).Value = 1452.3939
$ws_ARM.Cells.Item(122, 9).Value = 1282.7142
$ws_ARM.Cells.Item(122, 10).Value = 1749.3334
$ws_ARM.Cells.Item(122, 11).Value = 3848.1426
$ws_ARM.Cells.Item(122, 12).Value = 5248.0002
$ws_ARM.Cells.Item(122, 13).Value = -1398.1426
$ws_ARM.Cells.Item(122, 14).Value = -10148.0002
$ws_ARM.Cells.Item(132, 8).Value = 2648.4443
$ws_ARM.Cells.Item(132, 9).Value = 2234.5557
$ws_ARM.Cells.Item(132, 10).Value = 3476.2222
$ws_ARM.Cells.Item(132, 11).Value = 6703.6671
$ws_ARM.Cells.Item(132, 12).Value = 10428.6666
$ws_ARM.Cells.Item(132, 13).Value = -4173.6671
$ws_ARM.Cells.Item(132, 14).Value = -15488.6666
$ws_ARM.Cells.Item(136, 8).Value = 47620076
$ws_ARM.Cells.Item(136, 9).Value = 55556444
$ws_ARM.Cells.Item(136, 10).Value = 1853.3334
$ws_ARM.Cells.Item(136, 11).Value = 166669332
$ws_ARM.Cells.Item(136, 12).Value = 5560.0002
$ws_ARM.Cells.Item(136, 13).Value = -166666782
$ws_ARM.Cells.Item(136, 14).Value = -10660.0002
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(22, 8).Value = 250
$ws_BSM.Cells.Item(22, 9).Value = 350
$ws_BSM.Cells.Item(22, 10).Value = 50
$ws_BSM.Cells.Item(22, 11).Value = 350
$ws_BSM.Cells.Item(22, 12).Value = 50
$ws_BSM.Cells.Item(22, 13).Value = -177
$ws_BSM.Cells.Item(22, 14).Value = -396
$ws_BSM.Cells.Item(134, 8).Value = 1555.1052
$ws_BSM.Cells.Item(134, 9).Value = 1032.1765
$ws_BSM.Cells.Item(134, 10).Value = 6000
$ws_BSM.Cells.Item(134, 11).Value = 3096.5295
$ws_BSM.Cells.Item(134, 12).Value = 18000
$ws_BSM.Cells.Item(134, 13).Value = -561.5295000000001
$ws_BSM.Cells.Item(134, 14).Value = -23070
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(58, 8).Value = 1700.15
$ws_CRP.Cells.Item(58, 9).Value = 1379.3636
$ws_CRP.Cells.Item(58, 10).Value = 2092.2222
$ws_CRP.Cells.Item(58, 11).Value = 1379.3636
$ws_CRP.Cells.Item(58, 12).Value = 2092.2222
$ws_CRP.Cells.Item(58, 13).Value = -1176.3636
$ws_CRP.Cells.Item(58, 14).Value = -2498.2222
$ws_CRP.Cells.Item(74, 8).Value = 29166.666
$ws_CRP.Cells.Item(74, 9).Value = 22500
$ws_CRP.Cells.Item(74, 11).Value = 22500
$ws_CRP.Cells.Item(74, 13).Value = -21626
$ws_CRP.Cells.Item(77, 8).Value = 29166.666
$ws_CRP.Cells.Item(77, 9).Value = 22500
$ws_CRP.Cells.Item(77, 11).Value = 67500
$ws_CRP.Cells.Item(77, 13).Value = -63132
$ws_CRP.Cells.Item(136, 8).Value = 1700.15
$ws_CRP.Cells.Item(136, 9).Value = 1379.3636
$ws_CRP.Cells.Item(136, 10).Value = 2092.2222
$ws_CRP.Cells.Item(136, 11).Value = 4138.0908
$ws_CRP.Cells.Item(136, 12).Value = 6276.6666
$ws_CRP.Cells.Item(136, 13).Value = -1588.0908
$ws_CRP.Cells.Item(136, 14).Value = -11376.6666
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(107, 8).Value = 5289.6
$ws_CUL.Cells.Item(107, 9).Value = 377.25
$ws_CUL.Cells.Item(107, 10).Value = 6517.6875
$ws_CUL.Cells.Item(107, 11).Value = 1131.75
$ws_CUL.Cells.Item(107, 12).Value = 19553.0625
$ws_CUL.Cells.Item(107, 13).Value = 788.25
$ws_CUL.Cells.Item(107, 14).Value = -23393.0625
$ws_CUL.Cells.Item(124, 8).Value = 1487.1666
$ws_CUL.Cells.Item(124, 9).Value = 0
$ws_CUL.Cells.Item(124, 10).Value = 1487.1666
$ws_CUL.Cells.Item(124, 11).Value = 0
$ws_CUL.Cells.Item(124, 12).Value = 4461.4998
$ws_CUL.Cells.Item(124, 13).Value = $null
$ws_CUL.Cells.Item(124, 14).Value = -14281.4998
$ws_CUL.Cells.Item(131, 8).Value = 40002028
$ws_CUL.Cells.Item(131, 9).Value = 333333600
$ws_CUL.Cells.Item(131, 10).Value = 2267.6365
$ws_CUL.Cells.Item(131, 11).Value = 1000000800
$ws_CUL.Cells.Item(131, 12).Value = 6802.9095
$ws_CUL.Cells.Item(131, 13).Value = -999995760
$ws_CUL.Cells.Item(131, 14).Value = -16882.9095
$ws_CUL.Cells.Item(132, 8).Value = 921.5714
$ws_CUL.Cells.Item(132, 9).Value = 887.6875
$ws_CUL.Cells.Item(132, 11).Value = 7989.1875
$ws_CUL.Cells.Item(132, 13).Value = -5459.1875
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(102, 8).Value = 934.0526
$ws_GSM.Cells.Item(102, 9).Value = 916.7692
$ws_GSM.Cells.Item(102, 10).Value = 971.5
$ws_GSM.Cells.Item(102, 11).Value = 916.7692
$ws_GSM.Cells.Item(102, 12).Value = 971.5
$ws_GSM.Cells.Item(102, 13).Value = 705.2308
$ws_GSM.Cells.Item(102, 14).Value = -4215.5
$ws_GSM.Cells.Item(132, 8).Value = 5467.6665
$ws_GSM.Cells.Item(132, 9).Value = 6742.2
$ws_GSM.Cells.Item(132, 10).Value = 3874.5
$ws_GSM.Cells.Item(132, 11).Value = 20226.6
$ws_GSM.Cells.Item(132, 12).Value = 11623.5
$ws_GSM.Cells.Item(132, 13).Value = -17696.6
$ws_GSM.Cells.Item(132, 14).Value = -16683.5
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(22, 8).Value = 1075.1666
$ws_LTW.Cells.Item(22, 10).Value = 1234
$ws_LTW.Cells.Item(22, 12).Value = 1234
$ws_LTW.Cells.Item(22, 14).Value = -1824
$ws_LTW.Cells.Item(27, 8).Value = 1075.1666
$ws_LTW.Cells.Item(27, 10).Value = 1234
$ws_LTW.Cells.Item(27, 12).Value = 1234
$ws_LTW.Cells.Item(27, 14).Value = -1448
$ws_LTW.Cells.Item(40, 8).Value = 3229.318
$ws_LTW.Cells.Item(40, 9).Value = 2188.75
$ws_LTW.Cells.Item(40, 11).Value = 2188.75
$ws_LTW.Cells.Item(40, 13).Value = -2052.75
$ws_LTW.Cells.Item(82, 8).Value = 1380.3
$ws_LTW.Cells.Item(82, 10).Value = 1801.5
$ws_LTW.Cells.Item(82, 12).Value = 1801.5
$ws_LTW.Cells.Item(82, 14).Value = -2523.5
$ws_LTW.Cells.Item(85, 8).Value = 1380.3
$ws_LTW.Cells.Item(85, 10).Value = 1801.5
$ws_LTW.Cells.Item(85, 12).Value = 1801.5
$ws_LTW.Cells.Item(85, 14).Value = -4297.5
$ws_LTW.Cells.Item(100, 8).Value = 1218
$ws_LTW.Cells.Item(100, 9).Value = 1096.6666
$ws_LTW.Cells.Item(100, 11).Value = 1096.6666
$ws_LTW.Cells.Item(100, 13).Value = -555.6666
$ws_LTW.Cells.Item(122, 8).Value = 62501600
$ws_LTW.Cells.Item(122, 9).Value = 83334630
$ws_LTW.Cells.Item(122, 11).Value = 250003890
$ws_LTW.Cells.Item(122, 13).Value = -250001440
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(100, 8).Value = 206.63637
$ws_WVR.Cells.Item(100, 9).Value = 199.625
$ws_WVR.Cells.Item(100, 10).Value = 225.33333
$ws_WVR.Cells.Item(100, 11).Value = 399.25
$ws_WVR.Cells.Item(100, 12).Value = 450.66666
$ws_WVR.Cells.Item(100, 13).Value = 141.75
$ws_WVR.Cells.Item(100, 14).Value = -1532.66666
$ws_WVR.Cells.Item(132, 8).Value = 1378.3704
$ws_WVR.Cells.Item(132, 9).Value = 957.63635
$ws_WVR.Cells.Item(132, 10).Value = 3229.6
$ws_WVR.Cells.Item(132, 11).Value = 2872.90905
$ws_WVR.Cells.Item(132, 12).Value = 9688.799999999999
$ws_WVR.Cells.Item(132, 13).Value = -342.9090500000002
$ws_WVR.Cells.Item(132, 14).Value = -14748.8
$ws_WVR.Cells.Item(136, 8).Value = 991.64703
$ws_WVR.Cells.Item(136, 9).Value = 919.3077
$ws_WVR.Cells.Item(136, 11).Value = 2757.9231
$ws_WVR.Cells.Item(136, 13).Value = -207.9231
